$wb = $excel.ActiveWorkbook
$wsItem = $wb.Worksheets.Item("Item")

# Swap columns U (21) and V (22): cut V, insert before U.
# This moves V's data/width/style to column U's old position and shifts
# the old U column (and everything to its right) one place to the right,
# landing U's original content in column V.
$wsItem.Columns.Item(22).Cut()
$wsItem.Columns.Item(21).Insert()

# Make "Item" the active sheet/tab and move the selection to M7
# (this also clears the old topLeftCell="H1" scroll position and
# removes tabSelected from whichever sheet had it before).
$wsItem.Activate()
$wsItem.Range("M7").Select()
